$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column C (Quantity), shifting Quantity -> D and Cost Basis -> E
$ws.Columns.Item(3).Insert()

# New header for inserted column
$ws.Range("C1").Value = "Price"

# Populate Price values (Cost Basis / Quantity from original, higher-precision source data)
$ws.Range("C2").Value = 0.5448788169564495
$ws.Range("C3").Value = 0.06803486902880881
$ws.Range("C4").Value = 0.07333113036895253
$ws.Range("C5").Value = 0.08702114026711345
$ws.Range("C6").Value = 0.06583444599568089
$ws.Range("C7").Value = 0.06109126144076911
$ws.Range("C8").Value = 0.06220576971199036
$ws.Range("C9").Value = 1.00028758268002
$ws.Range("C10").Value = 0.05164678344993132
$ws.Range("C11").Value = 0.04581990110018071
$ws.Range("C12").Value = 0.05645797071004285
$ws.Range("C13").Value = 0.0584899292712627
$ws.Range("C14").Value = 0.704848484832252
$ws.Range("C15").Value = 0.7312994935923857
$ws.Range("C16").Value = 0.7241962123126896
$ws.Range("C17").Value = 0.7469625014359766
$ws.Range("C18").Value = 0.0484474033023344
$ws.Range("C19").Value = 0.05301888793675442
$ws.Range("C20").Value = 0.0513950184195149
$ws.Range("C21").Value = 0.428976624570404
$ws.Range("C22").Value = 0.4503605885935721
$ws.Range("C23").Value = 0.05890923103526676
$ws.Range("C24").Value = 0.4801375295817313
$ws.Range("C25").Value = 0.4488356714405042
$ws.Range("C26").Value = 0.4546727091699747
$ws.Range("C27").Value = 0.4371946406745632
$ws.Range("C28").Value = 0.06300502512508706
$ws.Range("C29").Value = 0.06458291457213151
$ws.Range("C30").Value = 0.05924242424243397
$ws.Range("C31").Value = 0.06411055276443217
$ws.Range("C32").Value = 0.0690242930183286
$ws.Range("C33").Value = 0.07290452261403836
$ws.Range("C34").Value = 0.07010918921294899
$ws.Range("C35").Value = 0.06417085427117038
$ws.Range("C36").Value = 0.06974874371841119
$ws.Range("C37").Value = 0.07608365565557711
$ws.Range("C38").Value = 0.05367839196031339
$ws.Range("C39").Value = 0.04459995351554337
$ws.Range("C40").Value = 0.04459144923770082
$ws.Range("C41").Value = 0.04455868101211714
$ws.Range("C42").Value = 0.04513567839200289
$ws.Range("C43").Value = 1178.442052468954
$ws.Range("C44").Value = 0.05807072647076265
$ws.Range("C45").Value = 0.05805059388054373
